$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values for rows 2..39, computed from the regenerated
# save_data (K replaces the old Strike# derived value).
$kValues = @(
    0, 1, 1, 1, 1, 3, 1, 1, 1, 0,
    2, 1, 0, 2, 1, 2, 1, 2, 0, 1,
    4, 0, 1, 1, 1, 0, 1, 1, 1, 0,
    3, 3, 0, 2, 2, 2, 0, 2
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
